$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.UnMerge()
$ws.Cells.Validation.Delete()
$ws.Cells.Clear()

$noteLead = "Dependent Dropdown list.    "
$noteTail = "Prequisite:- Name Range is required for dependent dropdown list."
$ws.Range("A2").Value = "$noteLead$noteTail"
$ws.Range("A2").Characters($noteLead.Length + 1, $noteTail.Length).Font.Size = 11

$ws.Range("A2:E6").Merge()

$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").WrapText = $true

"B2","C2","D2","E2","A3","B3","C3","D3","E3","A4","B4","C4","D4","E4","A5","B5","C5","D5","E5","A6","B6","C6","D6","E6" | ForEach-Object {
    $ws.Range($_).HorizontalAlignment = -4108
    $ws.Range($_).WrapText = $true
}
